# Iteration 14 branch compilation issues fixes
#
# Adds a new "Prepaid_Application_Upload" worksheet (a near-duplicate of the
# existing "LoyaltyPlan" sheet, with a new test-case name in A2), makes it
# the active/selected sheet, and tidies up the selection state that is left
# behind on "LoyaltyPlan".

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "LoyaltyPlan" to create the new sheet -------------------
$source = $wb.Worksheets.Item("LoyaltyPlan")
$source.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Prepaid_Application_Upload"

# --- 2. Update the new sheet's unique content ------------------------------
$newSheet.Range("A2").Value = "TC_Application_Upload_Prepaid"

# --- 3. Fix up view/selection state ----------------------------------------
# New sheet becomes the active tab, with A2 selected (mirrors a freshly
# duplicated sheet landing on its first data row).
$newSheet.Activate()
$newSheet.Range("A2").Select()

# Old "LoyaltyPlan" sheet no longer is the active tab; leave its selection
# as a full-used-range selection (A1:XFD2) instead of the old EB6 cell.
$loyaltyPlan = $wb.Worksheets.Item("LoyaltyPlan")
$loyaltyPlan.Range("A1:XFD2").Select()

# Re-activate the new sheet so it is the workbook's active/visible tab.
$newSheet.Activate()
